$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 54, shifting existing data (old rows 54-138)
# down to rows 57-141.
$ws.Rows("54:56").Insert()

# Populate the 3 newly inserted rows (54, 55, 56) with new data.
# Columns: A=Mercado ID, B=Mercado, C=Region, D=Fecha, E=Codreg, F=Tipo,
# G=Producto ID, H=Producto, I=Categoria ID, J=Categoria, K=Variedad,
# L=Calidad, M=Volumen, N=Precio minimo, O=Precio maximo,
# P=Precio promedio ponderado, Q=Unidad de comercializacion, R=Origen,
# S=Precio $/Kg, T=Kg / unidad

$commonA = 10
$commonB = "Vega Modelo de Temuco"
$commonC = "La Araucanía"
$commonE = 9
$commonF = "Fruta"
$commonG = 100101
$commonH = "Berries"
$commonI = 100112025
$commonJ = "Frutilla"
$commonK = "Sin especificar"
$commonT = 7

# Row 54
$ws.Cells.Item(54,1).Value = $commonA
$ws.Cells.Item(54,2).Value = $commonB
$ws.Cells.Item(54,3).Value = $commonC
$ws.Cells.Item(54,4).Value = 44477
$ws.Cells.Item(54,5).Value = $commonE
$ws.Cells.Item(54,6).Value = $commonF
$ws.Cells.Item(54,7).Value = $commonG
$ws.Cells.Item(54,8).Value = $commonH
$ws.Cells.Item(54,9).Value = $commonI
$ws.Cells.Item(54,10).Value = $commonJ
$ws.Cells.Item(54,11).Value = $commonK
$ws.Cells.Item(54,12).Value = "Primera"
$ws.Cells.Item(54,13).Value = 800
$ws.Cells.Item(54,14).Value = 14000
$ws.Cells.Item(54,15).Value = 14000
$ws.Cells.Item(54,16).Value = 14000
$ws.Cells.Item(54,17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(54,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(54,19).Value = 2000
$ws.Cells.Item(54,20).Value = $commonT

# Row 55
$ws.Cells.Item(55,1).Value = $commonA
$ws.Cells.Item(55,2).Value = $commonB
$ws.Cells.Item(55,3).Value = $commonC
$ws.Cells.Item(55,4).Value = 44477
$ws.Cells.Item(55,5).Value = $commonE
$ws.Cells.Item(55,6).Value = $commonF
$ws.Cells.Item(55,7).Value = $commonG
$ws.Cells.Item(55,8).Value = $commonH
$ws.Cells.Item(55,9).Value = $commonI
$ws.Cells.Item(55,10).Value = $commonJ
$ws.Cells.Item(55,11).Value = $commonK
$ws.Cells.Item(55,12).Value = "Segunda"
$ws.Cells.Item(55,13).Value = 200
$ws.Cells.Item(55,14).Value = 10000
$ws.Cells.Item(55,15).Value = 12000
$ws.Cells.Item(55,16).Value = 11000
$ws.Cells.Item(55,17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(55,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(55,19).Value = 1571
$ws.Cells.Item(55,20).Value = $commonT

# Row 56
$ws.Cells.Item(56,1).Value = $commonA
$ws.Cells.Item(56,2).Value = $commonB
$ws.Cells.Item(56,3).Value = $commonC
$ws.Cells.Item(56,4).Value = 44477
$ws.Cells.Item(56,5).Value = $commonE
$ws.Cells.Item(56,6).Value = $commonF
$ws.Cells.Item(56,7).Value = $commonG
$ws.Cells.Item(56,8).Value = $commonH
$ws.Cells.Item(56,9).Value = $commonI
$ws.Cells.Item(56,10).Value = $commonJ
$ws.Cells.Item(56,11).Value = $commonK
$ws.Cells.Item(56,12).Value = "Tercera"
$ws.Cells.Item(56,13).Value = 80
$ws.Cells.Item(56,14).Value = 7000
$ws.Cells.Item(56,15).Value = 7000
$ws.Cells.Item(56,16).Value = 7000
$ws.Cells.Item(56,17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(56,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(56,19).Value = 1000
$ws.Cells.Item(56,20).Value = $commonT
